$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New notes about installations (bootstrap, react-cookie, react-router-dom,
# reactstrap) and updates to the Home.js / index.js routing work.
$ws.Range("A63").Value = "create Home.js"
$ws.Range("A64").Value = "replace index.js to Home"

$ws.Range("A65").Value = "npm install bootstrap "
$ws.Range("D65").Value = "react-cookie"
$ws.Range("F65").Value = "react-router-dom AND reactstrap"

$ws.Range("A67").Value = "react-router routes between pages"

# Move the visible selection down to the newly added notes, same as the
# author scrolling the sheet to row 63 before saving.
$ws.Range("M63").Select()
